$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 ("time_taken"), matching the style of the existing
# header cells (bold font, centered, thin border) by copying the format
# from the neighboring header cell E1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row
$timeTaken = @{
    2 = "2021-10-05 13:41:12.251127"
    3 = "2021-10-05 13:41:12.251140"
    4 = "2021-10-05 13:41:12.251144"
    5 = "2021-10-05 13:41:12.251147"
    6 = "2021-10-05 13:41:12.251150"
    7 = "2021-10-05 13:41:12.251153"
    8 = "2021-10-05 13:41:12.251156"
    9 = "2021-10-05 13:41:12.251159"
    10 = "2021-10-05 13:41:12.251163"
    11 = "2021-10-05 13:41:12.251166"
    12 = "2021-10-05 13:41:12.251169"
    13 = "2021-10-05 13:41:12.251172"
    14 = "2021-10-05 13:41:12.251175"
    15 = "2021-10-05 13:41:12.251178"
    16 = "2021-10-05 13:41:12.251181"
    17 = "2021-10-05 13:41:12.251184"
    18 = "2021-10-05 13:41:12.251187"
    19 = "2021-10-05 13:41:12.251190"
    20 = "2021-10-05 13:41:12.251193"
    21 = "2021-10-05 13:41:12.251196"
    22 = "2021-10-05 13:41:12.251199"
    23 = "2021-10-05 13:41:12.251202"
    24 = "2021-10-05 13:41:12.251205"
    25 = "2021-10-05 13:41:12.251208"
    26 = "2021-10-05 13:41:12.251211"
    27 = "2021-10-05 13:41:12.251214"
    28 = "2021-10-05 13:41:12.251217"
    29 = "2021-10-05 13:41:12.251220"
    30 = "2021-10-05 13:41:12.251223"
    31 = "2021-10-05 13:41:12.251226"
    32 = "2021-10-05 13:41:12.251229"
    33 = "2021-10-05 13:41:12.251232"
    34 = "2021-10-05 13:41:12.251236"
    35 = "2021-10-05 13:41:12.251239"
    36 = "2021-10-05 13:41:12.251242"
    37 = "2021-10-05 13:41:12.251245"
    38 = "2021-10-05 13:41:12.251248"
    39 = "2021-10-05 13:41:12.251250"
    40 = "2021-10-05 13:41:12.251253"
    41 = "2021-10-05 13:41:12.251256"
    42 = "2021-10-05 13:41:12.251260"
    43 = "2021-10-05 13:41:12.251263"
    44 = "2021-10-05 13:41:12.251266"
    45 = "2021-10-05 13:41:12.251269"
    46 = "2021-10-05 13:41:12.251272"
    47 = "2021-10-05 13:41:12.251275"
    48 = "2021-10-05 13:41:12.251278"
    49 = "2021-10-05 13:41:12.251281"
    50 = "2021-10-05 13:41:12.251284"
    51 = "2021-10-05 13:41:12.251287"
    52 = "2021-10-05 13:41:12.251290"
    53 = "2021-10-05 13:41:12.251293"
    54 = "2021-10-05 13:41:12.251296"
    55 = "2021-10-05 13:41:12.251299"
    56 = "2021-10-05 13:41:12.251302"
    57 = "2021-10-05 13:41:12.251305"
    58 = "2021-10-05 13:41:12.251308"
    59 = "2021-10-05 13:41:12.251311"
    60 = "2021-10-05 13:41:12.251314"
    61 = "2021-10-05 13:41:12.251317"
    62 = "2021-10-05 13:41:12.251320"
    63 = "2021-10-05 13:41:12.251323"
    64 = "2021-10-05 13:41:12.251326"
    65 = "2021-10-05 13:41:12.251329"
    66 = "2021-10-05 13:41:12.251333"
}

foreach ($row in $timeTaken.Keys) {
    $ws.Cells.Item($row, 6).Value = $timeTaken[$row]
}

Write-Host "Done"
